# Auto-generated: update Leve profit calculation figures across sheets
# per scheduled runner refresh (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4946
$ws.Range("I51").Value = 4897.5
$ws.Range("K51").Value = 4897.5
$ws.Range("M51").Value = -4413.5
$ws.Range("H70").Value = 3476.9285
$ws.Range("J70").Value = 3607
$ws.Range("L70").Value = 10821
$ws.Range("N70").Value = -11361
$ws.Range("H73").Value = 3476.9285
$ws.Range("J73").Value = 3607
$ws.Range("L73").Value = 10821
$ws.Range("N73").Value = -12693
$ws.Range("H80").Value = 6279.4116
$ws.Range("I80").Value = 4245.615
$ws.Range("J80").Value = 12889.25
$ws.Range("K80").Value = 12736.845
$ws.Range("L80").Value = 38667.75
$ws.Range("M80").Value = -11738.845
$ws.Range("N80").Value = -40663.75
$ws.Range("H83").Value = 6279.4116
$ws.Range("I83").Value = 4245.615
$ws.Range("J83").Value = 12889.25
$ws.Range("K83").Value = 38210.535
$ws.Range("L83").Value = 116003.25
$ws.Range("M83").Value = -33218.535
$ws.Range("N83").Value = -125987.25
$ws.Range("H98").Value = 4096.196
$ws.Range("I98").Value = 4076.6667
$ws.Range("K98").Value = 4076.6667
$ws.Range("M98").Value = -2578.6667
$ws.Range("H106").Value = 3585.25
$ws.Range("I106").Value = 3463.6667
$ws.Range("J106").Value = 3950
$ws.Range("K106").Value = 3463.6667
$ws.Range("L106").Value = 3950
$ws.Range("M106").Value = -2832.6667
$ws.Range("N106").Value = -5212
$ws.Range("H107").Value = 628.1111
$ws.Range("J107").Value = 1002.5
$ws.Range("L107").Value = 1002.5
$ws.Range("N107").Value = -4842.5
$ws.Range("H122").Value = 4096.196
$ws.Range("I122").Value = 4076.6667
$ws.Range("K122").Value = 12230.0001
$ws.Range("M122").Value = -9780.000100000001
$ws.Range("H138").Value = 2786.4082
$ws.Range("J138").Value = 3055.303
$ws.Range("L138").Value = 9165.909
$ws.Range("N138").Value = -19445.909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4994.3887
$ws.Range("I102").Value = 3593.3333
$ws.Range("K102").Value = 3593.3333
$ws.Range("M102").Value = -1971.3333
$ws.Range("H122").Value = 2266
$ws.Range("I122").Value = 2266
$ws.Range("K122").Value = 6798
$ws.Range("M122").Value = -4348
$ws.Range("H132").Value = 3412.5557
$ws.Range("I132").Value = 2732.4666
$ws.Range("K132").Value = 8197.399800000001
$ws.Range("M132").Value = -5667.399800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5061.2905
$ws.Range("J99").Value = 7204.5
$ws.Range("L99").Value = 7204.5
$ws.Range("N99").Value = -10200.5
$ws.Range("H130").Value = 93181.82000000001
$ws.Range("J130").Value = 100000
$ws.Range("L130").Value = 100000
$ws.Range("N130").Value = -110040
$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5195.75
$ws.Range("I22").Value = 5195.75
$ws.Range("K22").Value = 5195.75
$ws.Range("M22").Value = -4845.75
$ws.Range("H99").Value = 5478.6
$ws.Range("I99").Value = 3118
$ws.Range("K99").Value = 3118
$ws.Range("M99").Value = -1620
$ws.Range("H126").Value = 5478.6
$ws.Range("I126").Value = 3118
$ws.Range("K126").Value = 9354
$ws.Range("M126").Value = -6884
$ws.Range("H139").Value = 139779.58
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 139779.58
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 139779.58
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -150059.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 294.55554
$ws.Range("I14").Value = 294.55554
$ws.Range("K14").Value = 883.66662
$ws.Range("M14").Value = -710.66662
$ws.Range("H60").Value = 999.8
$ws.Range("I60").Value = 999.8
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 2999.4
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -2748.4
$ws.Range("N60").ClearContents()
$ws.Range("H86").Value = 509.66666
$ws.Range("J86").Value = 557.1667
$ws.Range("L86").Value = 1671.5001
$ws.Range("N86").Value = -4043.5001
$ws.Range("H89").Value = 509.66666
$ws.Range("J89").Value = 557.1667
$ws.Range("L89").Value = 5014.5003
$ws.Range("N89").Value = -16870.5003
$ws.Range("H132").Value = 1296.1538
$ws.Range("J132").Value = 1341.6666
$ws.Range("L132").Value = 12074.9994
$ws.Range("N132").Value = -17134.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3461.3333
$ws.Range("I80").Value = 3274.4443
$ws.Range("J80").Value = 4022
$ws.Range("K80").Value = 3274.4443
$ws.Range("L80").Value = 4022
$ws.Range("M80").Value = -2276.4443
$ws.Range("N80").Value = -6018
$ws.Range("H83").Value = 3461.3333
$ws.Range("I83").Value = 3274.4443
$ws.Range("J83").Value = 4022
$ws.Range("K83").Value = 16372.2215
$ws.Range("L83").Value = 20110
$ws.Range("M83").Value = -11380.2215
$ws.Range("N83").Value = -30094
$ws.Range("H93").Value = 36324.383
$ws.Range("J93").Value = 36324.383
$ws.Range("L93").Value = 36324.383
$ws.Range("N93").Value = -40068.383
$ws.Range("H97").Value = 2010.2
$ws.Range("I97").Value = 518.1111
$ws.Range("K97").Value = 518.1111
$ws.Range("M97").Value = -22.11109999999996
$ws.Range("H122").Value = 2124.1304
$ws.Range("J122").Value = 1955.9
$ws.Range("L122").Value = 5867.700000000001
$ws.Range("N122").Value = -10767.7
$ws.Range("H126").Value = 1749.3077
$ws.Range("I126").Value = 1496.875
$ws.Range("J126").Value = 2153.2
$ws.Range("K126").Value = 4490.625
$ws.Range("L126").Value = 6459.599999999999
$ws.Range("M126").Value = -2020.625
$ws.Range("N126").Value = -11399.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 552.625
$ws.Range("I22").Value = 653.75
$ws.Range("K22").Value = 653.75
$ws.Range("M22").Value = -358.75
$ws.Range("H27").Value = 552.625
$ws.Range("I27").Value = 653.75
$ws.Range("K27").Value = 653.75
$ws.Range("M27").Value = -546.75
$ws.Range("H55").Value = 646.4583
$ws.Range("I55").Value = 806.53845
$ws.Range("J55").Value = 457.27274
$ws.Range("K55").Value = 806.53845
$ws.Range("L55").Value = 457.27274
$ws.Range("M55").Value = -633.53845
$ws.Range("N55").Value = -803.27274
$ws.Range("H61").Value = 78705.92
$ws.Range("I61").Value = 84764.75
$ws.Range("K61").Value = 84764.75
$ws.Range("M61").Value = -84562.75
$ws.Range("H68").Value = 10980.223
$ws.Range("I68").Value = 10332.143
$ws.Range("J68").Value = 13248.5
$ws.Range("K68").Value = 10332.143
$ws.Range("L68").Value = 13248.5
$ws.Range("M68").Value = -9583.143
$ws.Range("N68").Value = -14746.5
$ws.Range("H71").Value = 10980.223
$ws.Range("I71").Value = 10332.143
$ws.Range("J71").Value = 13248.5
$ws.Range("K71").Value = 51660.715
$ws.Range("L71").Value = 66242.5
$ws.Range("M71").Value = -47916.715
$ws.Range("N71").Value = -73730.5
$ws.Range("H100").Value = 2384985.8
$ws.Range("I100").Value = 4169975.2
$ws.Range("J100").Value = 4999.6665
$ws.Range("K100").Value = 4169975.2
$ws.Range("L100").Value = 4999.6665
$ws.Range("M100").Value = -4169434.2
$ws.Range("N100").Value = -6081.6665
$ws.Range("H113").Value = 78705.92
$ws.Range("I113").Value = 84764.75
$ws.Range("K113").Value = 84764.75
$ws.Range("M113").Value = -82594.75
$ws.Range("H122").Value = 4721.143
$ws.Range("I122").Value = 3377.2222
$ws.Range("J122").Value = 7140.2
$ws.Range("K122").Value = 10131.6666
$ws.Range("L122").Value = 21420.6
$ws.Range("M122").Value = -7681.6666
$ws.Range("N122").Value = -26320.6
$ws.Range("H132").Value = 6303.7104
$ws.Range("I132").Value = 6212.1562
$ws.Range("J132").Value = 6792
$ws.Range("K132").Value = 18636.4686
$ws.Range("L132").Value = 20376
$ws.Range("M132").Value = -16106.4686
$ws.Range("N132").Value = -25436
$ws.Range("H134").Value = 31000
$ws.Range("I134").Value = 31000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 31000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -25930
$ws.Range("N134").ClearContents()
